# Scheduled-runner profit/price refresh: pushes newly-pulled market-board
# data (H/I/J/K/L source values) into each sheet's Leve-profit table; the
# dependent M/N margin columns are literal numbers in this workbook (not
# formulas), so their refreshed results are written explicitly too.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 113.25
$ws.Range("I12").Value = 113.25
$ws.Range("K12").Value = 113.25
$ws.Range("M12").Value = 56.75
$ws.Range("H15").Value = 225865.42
$ws.Range("I15").Value = 225865.42
$ws.Range("K15").Value = 677596.26
$ws.Range("M15").Value = -677427.26
$ws.Range("H19").Value = 1103.1904
$ws.Range("I19").Value = 887.0909
$ws.Range("K19").Value = 887.0909
$ws.Range("M19").Value = -712.0909
$ws.Range("H33").Value = 292
$ws.Range("I33").Value = 204.7
$ws.Range("K33").Value = 204.7
$ws.Range("M33").Value = 24.30000000000001
$ws.Range("H76").Value = 200003200
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 200003200
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H96").Value = 125001180
$ws.Range("I96").Value = 1819.8
$ws.Range("K96").Value = 5459.4
$ws.Range("M96").Value = -4086.4
$ws.Range("H116").Value = 15876192
$ws.Range("I116").Value = 37039036
$ws.Range("J116").Value = 4059.5
$ws.Range("K116").Value = 37039036
$ws.Range("L116").Value = 4059.5
$ws.Range("M116").Value = -37035594
$ws.Range("N116").Value = -10943.5
$ws.Range("H137").Value = 7429.25
$ws.Range("I137").Value = 5114.857
$ws.Range("J137").Value = 11847.637
$ws.Range("K137").Value = 15344.571
$ws.Range("L137").Value = 35542.911
$ws.Range("M137").Value = -12794.571
$ws.Range("N137").Value = -40642.911

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4635057
$ws.Range("I32").Value = 2878.0715
$ws.Range("K32").Value = 2878.0715
$ws.Range("M32").Value = -2591.0715
$ws.Range("H45").Value = 2699.875
$ws.Range("I45").Value = 2519.8
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 2519.8
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = -2142.8
$ws.Range("N45").Value = -3754
$ws.Range("H61").Value = 6771.857
$ws.Range("I61").Value = 6715.846
$ws.Range("J61").Value = 7500
$ws.Range("K61").Value = 6715.846
$ws.Range("L61").Value = 7500
$ws.Range("M61").Value = -6503.846
$ws.Range("N61").Value = -7924
$ws.Range("H97").Value = 2123.1667
$ws.Range("I97").Value = 1499.8
$ws.Range("J97").Value = 2568.4285
$ws.Range("K97").Value = 1499.8
$ws.Range("L97").Value = 2568.4285
$ws.Range("M97").Value = -1003.8
$ws.Range("N97").Value = -3560.4285
$ws.Range("H122").Value = 5299.778
$ws.Range("I122").Value = 2599.5
$ws.Range("K122").Value = 7798.5
$ws.Range("M122").Value = -5348.5
$ws.Range("H132").Value = 951342.4
$ws.Range("I132").Value = 1124361
$ws.Range("J132").Value = 172758.67
$ws.Range("K132").Value = 3373083
$ws.Range("L132").Value = 518276.01
$ws.Range("M132").Value = -3370553
$ws.Range("N132").Value = -523336.01
$ws.Range("H136").Value = 6771.857
$ws.Range("I136").Value = 6715.846
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 20147.538
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = -17597.538
$ws.Range("N136").Value = -27600

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 12355504
$ws.Range("I80").Value = 1210.4166
$ws.Range("J80").Value = 22238938
$ws.Range("K80").Value = 1210.4166
$ws.Range("L80").Value = 22238938
$ws.Range("M80").Value = -212.4166
$ws.Range("N80").Value = -22240934
$ws.Range("H83").Value = 12355504
$ws.Range("I83").Value = 1210.4166
$ws.Range("J83").Value = 22238938
$ws.Range("K83").Value = 6052.083000000001
$ws.Range("L83").Value = 111194690
$ws.Range("M83").Value = -1060.083000000001
$ws.Range("N83").Value = -111204674

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 29566.572
$ws.Range("I41").Value = 21000
$ws.Range("J41").Value = 32993.2
$ws.Range("K41").Value = 21000
$ws.Range("L41").Value = 32993.2
$ws.Range("M41").Value = -20572
$ws.Range("N41").Value = -33849.2
$ws.Range("H132").Value = 3971.8965
$ws.Range("I132").Value = 3547.3818
$ws.Range("K132").Value = 10642.1454
$ws.Range("M132").Value = -8112.145400000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11000192
$ws.Range("I4").Value = 11000192
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 33000576
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = -33000464
$ws.Range("H63").Value = 13555
$ws.Range("H66").Value = 13555
$ws.Range("H68").Value = 715427.9
$ws.Range("I68").Value = 1332.5
$ws.Range("K68").Value = 3997.5
$ws.Range("M68").Value = -3186.5
$ws.Range("H71").Value = 715427.9
$ws.Range("I71").Value = 1332.5
$ws.Range("K71").Value = 11992.5
$ws.Range("M71").Value = -7936.5
$ws.Range("H75").Value = 3995
$ws.Range("J75").Value = 3995
$ws.Range("L75").Value = 11985
$ws.Range("N75").Value = -13981
$ws.Range("H78").Value = 3995
$ws.Range("J78").Value = 3995
$ws.Range("L78").Value = 35955
$ws.Range("N78").Value = -45939
$ws.Range("H107").Value = 4860.793
$ws.Range("I107").Value = 716.6667
$ws.Range("J107").Value = 5338.9614
$ws.Range("K107").Value = 2150.0001
$ws.Range("L107").Value = 16016.8842
$ws.Range("M107").Value = -230.0001000000002
$ws.Range("N107").Value = -19856.8842

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 58834950
$ws.Range("I80").Value = 142860100
$ws.Range("J80").Value = 17348.9
$ws.Range("K80").Value = 142860100
$ws.Range("L80").Value = 17348.9
$ws.Range("M80").Value = -142859102
$ws.Range("N80").Value = -19344.9
$ws.Range("H83").Value = 58834950
$ws.Range("I83").Value = 142860100
$ws.Range("J83").Value = 17348.9
$ws.Range("K83").Value = 714300500
$ws.Range("L83").Value = 86744.5
$ws.Range("M83").Value = -714295508
$ws.Range("N83").Value = -96728.5
$ws.Range("H97").Value = 994.3570999999999
$ws.Range("I97").Value = 912.0909
$ws.Range("K97").Value = 912.0909
$ws.Range("M97").Value = -416.0909
$ws.Range("H102").Value = 887134.0600000001
$ws.Range("I102").Value = 1555390.1
$ws.Range("J102").Value = 6251.136
$ws.Range("K102").Value = 1555390.1
$ws.Range("L102").Value = 6251.136
$ws.Range("M102").Value = -1553768.1
$ws.Range("N102").Value = -9495.136
$ws.Range("H126").Value = 71439400
$ws.Range("I126").Value = 250002400
$ws.Range("K126").Value = 750007200
$ws.Range("M126").Value = -750004730
$ws.Range("H132").Value = 34486444
$ws.Range("I132").Value = 40003068
$ws.Range("J132").Value = 7546.5
$ws.Range("K132").Value = 120009204
$ws.Range("L132").Value = 22639.5
$ws.Range("M132").Value = -120006674
$ws.Range("N132").Value = -27699.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 35716292
$ws.Range("I46").Value = 733
$ws.Range("J46").Value = 45456900
$ws.Range("K46").Value = 733
$ws.Range("L46").Value = 45456900
$ws.Range("M46").Value = -545
$ws.Range("N46").Value = -45457276
$ws.Range("H106").Value = 22999
$ws.Range("J106").Value = 22999
$ws.Range("L106").Value = 22999
$ws.Range("N106").Value = -25523
$ws.Range("H132").Value = 7212.7
$ws.Range("I132").Value = 6556.143
$ws.Range("K132").Value = 19668.429
$ws.Range("M132").Value = -17138.429
$ws.Range("H136").Value = 75006160
$ws.Range("I136").Value = 45461444
$ws.Range("K136").Value = 136384332
$ws.Range("M136").Value = -136381782

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 8999.5
$ws.Range("J105").Value = 8999.5
$ws.Range("L105").Value = 8999.5
$ws.Range("N105").Value = -15987.5
$ws.Range("H122").Value = 2976.4285
$ws.Range("I122").Value = 2977.8333
$ws.Range("J122").Value = 2968
$ws.Range("K122").Value = 8933.499899999999
$ws.Range("L122").Value = 8904
$ws.Range("M122").Value = -6483.499899999999
$ws.Range("N122").Value = -13804
$ws.Range("H126").Value = 7218.875
$ws.Range("I126").Value = 6760.25
$ws.Range("K126").Value = 20280.75
$ws.Range("M126").Value = -17810.75
